$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 34.33416056
$ws.Range("C2").Value = 134.01512139
$ws.Range("J2").Value = 0.41666666666666674
$ws.Range("K2").Value = 0.7083333333333333

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 34.33337361
$ws.Range("C3").Value = 134.01533389
$ws.Range("H3").ClearContents()
$ws.Range("J3").Value = 0.375
$ws.Range("K3").Value = 0.7083333333333333

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 34.25399139
$ws.Range("C4").Value = 134.029675
$ws.Range("H4").ClearContents()
$ws.Range("J4").Value = 0.375
$ws.Range("K4").Value = 0.7083333333333333

# --- Row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 34.30793111
$ws.Range("C5").Value = 134.10039139
$ws.Range("H5").ClearContents()
$ws.Range("J5").Value = 0.375
$ws.Range("K5").Value = 0.7083333333333333

# --- Row 6 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 34.17610333
$ws.Range("C6").Value = 134.07563889
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()

# --- Row 7 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 34.1375775
$ws.Range("C7").Value = 134.07610361
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()

# --- Row 8 ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 34.30207
$ws.Range("C8").Value = 133.96305
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()

# --- Row 9 ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 34.283235
$ws.Range("C9").Value = 133.960837
$ws.Range("E9").Value = "高松市国分寺町新名1432-2"
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("K9").ClearContents()

# --- Row 10 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 34.30739917
$ws.Range("C10").Value = 134.10019111
$ws.Range("H10").ClearContents()
$ws.Range("J10").Value = 0.375
$ws.Range("K10").Value = 0.7083333333333333
